# Update the "想去人数" (interest count) figures for three 合肥漫展 events.
# These numbers are refreshed by the site's scraper/build job (gh-pages
# "output generated at <commit>"), so only the numeric values change,
# on both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# "展览" (Exhibition) sheet - row numbers match directly.
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F5").Value = 5003
$wsExhibit.Range("F9").Value = 757
$wsExhibit.Range("F10").Value = 241

# "全部类型" (All types) sheet - same events, but the third one lives one
# row further down because this sheet merges entries from other sheets.
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F5").Value = 5003
$wsAll.Range("F9").Value = 757
$wsAll.Range("F11").Value = 241
